# Daily attendance processing - 2025-11-20 18:29:51
# Normalizes the "Recorded By" (column G) entries so that any list whose
# first item is not "System"/"system" gets that leading item rotated to
# the end of the list (i.e. "System" is moved to the front when present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ',\s*'

    if ($parts.Count -gt 1 -and $parts[0].Trim().ToLower() -ne 'system') {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $newVal = [string]::Join(', ', $rotated)
        $cell.Value2 = $newVal
    }
}
